# Apply cryptos.xlsx update (Sat Jun  8 10:52:18 UTC 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price-column values look like plain numbers (e.g. "7.39", "0.999").
# The source data stores these as text, so force a text number format on
# those specific cells before assigning, otherwise Excel's normal
# type-inference would convert them into numeric values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.405.39"
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "3.688.83"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "687.30"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "161.52"
$ws.Range("E6").Value = "  -5.82%  "
$ws.Range("D7").Value = "3.685.89"
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -5.54%  "
$ws.Range("E10").Value = "  -8.37%  "
$ws.Range("D11").Value = "7.39"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  -7.41%  "
$ws.Range("E13").Value = "  -5.97%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "33.09"
$ws.Range("E14").Value = "  -8.09%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.304.44"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("D16").Value = "3.678.91"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").Value = "69.410.25"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "16.07"
$ws.Range("E19").Value = "  -7.93%  "
$ws.Range("E20").Value = "  -9.09%  "
$ws.Range("D21").Value = "476.83"
$ws.Range("E21").Value = "  -7.34%  "
$ws.Range("E22").Value = "  -5.34%  "
$ws.Range("E23").Value = "  -7.80%  "
$ws.Range("D24").Value = "79.87"
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("D25").Value = "3.828.31"
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("E26").Value = "  -9.38%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "11.24"
$ws.Range("E28").Value = "  -8.06%  "
$ws.Range("D29").Value = "9.33"
$ws.Range("E29").Value = "  -9.69%  "
$ws.Range("E30").Value = "  -10.92%  "
$ws.Range("E31").Value = "  -10.26%  "
$ws.Range("E32").Value = "  -7.50%  "
$ws.Range("E33").Value = "  -7.84%  "
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "26.91"
$ws.Range("E36").Value = "  -7.48%  "
$ws.Range("D37").Value = "3.647.72"
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("D38").Value = "8.38"
$ws.Range("E38").Value = "  -8.30%  "
$ws.Range("D39").Value = "6.24"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").Value = "2.33"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("E41").Value = "  -8.91%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "0.952"
$ws.Range("E44").Value = "  -6.37%  "
$ws.Range("D45").Value = "163.64"
$ws.Range("E45").Value = "  -5.60%  "
$ws.Range("D46").Value = "48.17"
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("D47").Value = "29.89"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  -15.70%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "1.33"
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("D50").Value = "0.000281"
$ws.Range("E50").Value = "  -8.62%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "1.13"
$ws.Range("E51").Value = "  -1.55%  "
